$d = $word.ActiveDocument

function Set-NewParagraphText($range, $text, $underline) {
    $range.Font.Name = "Times New Roman"
    $range.Font.NameAscii = "Times New Roman"
    $range.Font.NameFarEast = "Times New Roman"
    $range.Font.NameOther = "Times New Roman"
    $range.Font.NameBi = "Times New Roman"
    $range.Font.Size = 12
    $range.Font.SizeBi = 12
    $range.Text = $text
    if ($underline) {
        $range.Font.Underline = 1
    }
}

# --- Insert the "Instructor Guide" heading as a brand new first paragraph ---
$firstPara = $d.Paragraphs.Item(1)
$firstPara.Range.InsertParagraphBefore()
$titlePara = $d.Paragraphs.Item(1)
Set-NewParagraphText $titlePara.Range "Instructor Guide" $true

# --- Append the "Common Pitfalls" section just before the trailing empty paragraph ---
$lastIndex = $d.Paragraphs.Count
$lastPara = $d.Paragraphs.Item($lastIndex)
$lastPara.Range.InsertParagraphBefore()

$lastIndex = $d.Paragraphs.Count
$lastPara = $d.Paragraphs.Item($lastIndex)
$lastPara.Range.InsertParagraphBefore()
$pitfallsTitlePara = $d.Paragraphs.Item($lastIndex)
Set-NewParagraphText $pitfallsTitlePara.Range "Common Pitfalls for Students and Instructors" $true

$lastIndex = $d.Paragraphs.Count
$lastPara = $d.Paragraphs.Item($lastIndex)
$lastPara.Range.InsertParagraphBefore()
$pitfallsBodyPara = $d.Paragraphs.Item($lastIndex)
Set-NewParagraphText $pitfallsBodyPara.Range "It is likely that students, especially those with no programming experience, may have some difficulty understanding the concept of shared memory parallelism, or even parallelism in general. Parallelism is a very different way of thinking than we normally think in our day to day lives, so it may take some students longer than expected to grasp the concept. Similarly, finding new and innovative ways for instructors to teach the same ideas might be difficult as well. This module attempts to provide several different examples to show different ways of teaching the same concept." $false
